$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Header text updates (volume number, report week dates)
# ---------------------------------------------------------------
$ws.Range("A8").Value = "Volume 30   Number  23"
$ws.Range("C9").Value = "Report Covering the Week  6/5/2023  Through  6/11/2023"

# ---------------------------------------------------------------
# Row 15 (Rape) - values only, no style/type changes
# ---------------------------------------------------------------
$ws.Range("I15").Value = 5
$ws.Range("K15").Value = 150
$ws.Range("M15").Value = -28.571428571428
$ws.Range("N15").Value = -58.333333333333

# ---------------------------------------------------------------
# Row 16 (Robbery)
# C16: numeric 1 -> text "0" (reuse style/shared-string pattern from C14)
# ---------------------------------------------------------------
$ws.Range("C14").Copy($ws.Range("C16"))
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = -100
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 6
$ws.Range("H16").Value = -83.333333333333
$ws.Range("J16").Value = 27
$ws.Range("K16").Value = -3.703703703703
$ws.Range("L16").Value = 52.941176470588
$ws.Range("M16").Value = -44.680851063829
$ws.Range("N16").Value = -80.303030303030

# ---------------------------------------------------------------
# Row 17 (Fel. Assault) - values only
# ---------------------------------------------------------------
$ws.Range("C17").Value = 1
$ws.Range("E17").Value = -66.666666666666
$ws.Range("F17").Value = 11
$ws.Range("G17").Value = 10
$ws.Range("H17").Value = 10
$ws.Range("I17").Value = 85
$ws.Range("J17").Value = 41
$ws.Range("K17").Value = 107.317073170732
$ws.Range("L17").Value = 54.545454545454
$ws.Range("M17").Value = 30.769230769230
$ws.Range("N17").Value = -27.350427350427

# ---------------------------------------------------------------
# Row 18 (Burglary)
# C18: numeric 2 -> text "0"
# ---------------------------------------------------------------
$ws.Range("C14").Copy($ws.Range("C18"))
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = -100
$ws.Range("F18").Value = 5
$ws.Range("G18").Value = 6
$ws.Range("H18").Value = -16.666666666666
$ws.Range("J18").Value = 31
$ws.Range("K18").Value = 32.258064516129
$ws.Range("L18").Value = 24.242424242424
$ws.Range("M18").Value = -48.75
$ws.Range("N18").Value = -93.166666666666

# ---------------------------------------------------------------
# Row 19 (Gr. Larceny) - values only
# ---------------------------------------------------------------
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = 57.142857142857
$ws.Range("F19").Value = 39
$ws.Range("G19").Value = 20
$ws.Range("H19").Value = 95
$ws.Range("I19").Value = 215
$ws.Range("J19").Value = 129
$ws.Range("K19").Value = 66.666666666666
$ws.Range("L19").Value = 50.349650349650
$ws.Range("M19").Value = 30.303030303030
$ws.Range("N19").Value = -43.421052631578

# ---------------------------------------------------------------
# Row 20 (G.L.A.)
# C20: text "0" -> numeric 2 (reuse numeric style from D20)
# ---------------------------------------------------------------
$ws.Range("D20").Copy($ws.Range("C20"))
$ws.Range("C20").Value = 2
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 7
$ws.Range("G20").Value = 12
$ws.Range("H20").Value = -41.666666666666
$ws.Range("I20").Value = 42
$ws.Range("J20").Value = 53
$ws.Range("K20").Value = -20.754716981132
$ws.Range("L20").Value = 90.909090909090
$ws.Range("M20").Value = 23.529411764705
$ws.Range("N20").Value = -96.579804560260

# ---------------------------------------------------------------
# Row 21 (TOTAL) - values only
# ---------------------------------------------------------------
$ws.Range("C21").Value = 15
$ws.Range("D21").Value = 15
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 65
$ws.Range("G21").Value = 55
$ws.Range("H21").Value = 18.181818181818
$ws.Range("I21").Value = 415
$ws.Range("J21").Value = 284
$ws.Range("K21").Value = 46.126760563380
$ws.Range("L21").Value = 50.362318840579
$ws.Range("M21").Value = 4.271356783919
$ws.Range("N21").Value = -83.198380566801

# ---------------------------------------------------------------
# Row 23 (Housing)
# D23: numeric 2 -> text "0" (reuse from C23, already text "0")
# E23: numeric -100 -> text "***.*" (reuse from E14)
# ---------------------------------------------------------------
$ws.Range("C23").Copy($ws.Range("D23"))
$ws.Range("E14").Copy($ws.Range("E23"))

# ---------------------------------------------------------------
# Row 24 (Petit Larceny) - values only
# ---------------------------------------------------------------
$ws.Range("C24").Value = 11
$ws.Range("D24").Value = 14
$ws.Range("E24").Value = -21.428571428571
$ws.Range("F24").Value = 90
$ws.Range("G24").Value = 54
$ws.Range("H24").Value = 66.666666666666
$ws.Range("I24").Value = 495
$ws.Range("J24").Value = 311
$ws.Range("K24").Value = 59.163987138263
$ws.Range("L24").Value = 136.842105263158
$ws.Range("M24").Value = -32.284541723666

# ---------------------------------------------------------------
# Row 25 (Misd. Assault) - values only
# ---------------------------------------------------------------
$ws.Range("C25").Value = 7
$ws.Range("D25").Value = 9
$ws.Range("E25").Value = -22.222222222222
$ws.Range("F25").Value = 31
$ws.Range("H25").Value = -20.512820512820
$ws.Range("I25").Value = 139
$ws.Range("J25").Value = 159
$ws.Range("K25").Value = -12.578616352201
$ws.Range("L25").Value = 31.132075471698
$ws.Range("M25").Value = -46.946564885496

# ---------------------------------------------------------------
# Row 26 (UCR Rape*)
# D26: text "0" -> numeric 1 (reuse numeric style from C26, while it still has value 1)
# E26: text "***.*" -> numeric 100 (reuse numeric style from H26)
# C26: numeric 1 -> numeric 2 (do after D26 copy so copy source is style-only correct)
# ---------------------------------------------------------------
$ws.Range("C26").Copy($ws.Range("D26"))
$ws.Range("D26").Value = 1
$ws.Range("H26").Copy($ws.Range("E26"))
$ws.Range("E26").Value = 100
$ws.Range("C26").Value = 2
$ws.Range("F26").Value = 5
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = 150
$ws.Range("I26").Value = 8
$ws.Range("J26").Value = 9
$ws.Range("K26").Value = -11.111111111111
$ws.Range("L26").Value = 14.285714285714

# ---------------------------------------------------------------
# Row 27 (Other Sex Crimes) - values only
# ---------------------------------------------------------------
$ws.Range("F27").Value = 5
$ws.Range("G27").Value = 5
$ws.Range("I27").Value = 16
$ws.Range("J27").Value = 19
$ws.Range("K27").Value = -15.789473684210
$ws.Range("L27").Value = 14.285714285714
